$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.914.64'
$ws.Range('D3').Value = '1.642.52'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = "'215.17"
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').Value = "'0.5040"
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').Value = "'1.005"
$ws.Range('E7').Value = '  -1.07%  '
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('D11').Value = "'0.07803"
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '1.653.64'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = '1.866.61'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').Value = "'0.5428"
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D16').Value = '0.0₅7862'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').Value = "'64.74"
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '25.953.31'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').Value = "'1.006"
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').Value = "'198.58"
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('D21').Value = "'4.388"
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').Value = "'9.945"
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = "'5.978"
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').Value = "'1.874"
$ws.Range('E25').Value = '  -4.63%  '
$ws.Range('D26').Value = "'140.11"
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('E27').Value = '  -1.74%  '
$ws.Range('D28').Value = "'6.853"
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').Value = "'1.242"
$ws.Range('E30').Value = '  -0.73%  '
$ws.Range('D31').Value = "'0.04898"
$ws.Range('E31').Value = '  -3.77%  '
$ws.Range('D32').Value = "'3.259"
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('D35').Value = "'2.368"
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').Value = "'0.8925"
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').Value = "'2.608"
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('D38').Value = '1.138.64'
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').Value = "'0.5550"
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('E40').Value = '  -1.34%  '
$ws.Range('D41').Value = "'1.004"
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('D42').Value = "'5.688"
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = "'0.8163"
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').Value = '0.0₈121'
$ws.Range('E45').Value = '  +5.57%  '
$ws.Range('D46').Value = '1.777.12'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = "'0.4535"
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('D49').Value = "'55.29"
$ws.Range('D50').Value = "'0.05089"
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = "'1.006"
$ws.Range('E51').Value = '  -0.88%  '
